# Applies the RcOnline/base.xlsx commit:
#   - Companies sheet: link row 39 (РИЦ-Облкоммунэнерго Ч/С (Алапаевский)) to new AccountTaskId 643
#   - Tasks sheet: update attempt/total/success/error counters and status/error messages
#     for several existing tasks, and append a brand-new task row (Id 643, Status 0)

$wb = $excel.ActiveWorkbook

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsTasks     = $wb.Worksheets.Item("Tasks")

# ------------------------------------------------------------------
# Companies sheet: new AccountTaskId for "РИЦ-Облкоммунэнерго Ч/С (Алапаевский)"
# ------------------------------------------------------------------
$wsCompanies.Range("D39").Value2 = 643

# ------------------------------------------------------------------
# Tasks sheet: per-row updates
# ------------------------------------------------------------------

# Row 128
$wsTasks.Range("E128").Value2 = 3468
$wsTasks.Range("F128").Value2 = 831

# Row 176
$wsTasks.Range("E176").Value2 = 702
$wsTasks.Range("F176").Value2 = 427
$wsTasks.Range("G176").Value2 = "Найдено 0 ПУ в ГИС ЖКХ (обл Свердловская,г Сухой Лог,ул Рябиновая,20)"

# Row 180
$wsTasks.Range("B180").Value2 = 3
$wsTasks.Range("E180").Value2 = 1273
$wsTasks.Range("F180").Value2 = 732
$wsTasks.Range("G180").Value2 = "Выгрузка начислений завершена."

# Row 198
$wsTasks.Range("D198").Value2 = 35226
$wsTasks.Range("E198").Value2 = 63
$wsTasks.Range("G198").Value2 = "ПД выгружен успешно"

# Row 199
$wsTasks.Range("B199").Value2 = 3
$wsTasks.Range("F199").Value2 = 4449
$wsTasks.Range("G199").Value2 = "Выгрузка начислений завершена."

# Row 200
$wsTasks.Range("B200").Value2 = 3
$wsTasks.Range("F200").Value2 = 4449
$wsTasks.Range("G200").Value2 = "Выгрузка начислений завершена."

# Row 202
$wsTasks.Range("B202").Value2 = 3
$wsTasks.Range("F202").Value2 = 2977
$wsTasks.Range("G202").Value2 = "Выгрузка начислений завершена."

# Row 206
$wsTasks.Range("D206").Value2 = 28279
$wsTasks.Range("E206").Value2 = 6
$wsTasks.Range("F206").Value2 = 5732
$wsTasks.Range("G206").Value2 = "Не найдены услуги в ГИС ЖКХ: 'Вывоз ЖБО, ТБО'"

# Row 207
$wsTasks.Range("B207").Value2 = 3
$wsTasks.Range("D207").Value2 = 11690
$wsTasks.Range("F207").Value2 = 11690
$wsTasks.Range("G207").Value2 = "Выгрузка начислений завершена."

# New row 209: freshly queued task, not yet processed
$wsTasks.Range("A209").Value2 = 643
$wsTasks.Range("B209").Value2 = 0
